$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation was inserted into the daily log at row 434. This
# pushes the former rows 434:490 down to 435:491 (dimension grows from
# A1:R490 to A1:R491) and fills the vacated row 434 with the new record.
$ws.Rows(434).Insert()

$ws.Range("A434").Value = 10
$ws.Range("B434").Value = "Vega Modelo de Temuco"
$ws.Range("C434").Value = "La Araucanía"
$ws.Range("D434").Value = [DateTime]::FromOADate(45142)
$ws.Range("E434").Value = 9
$ws.Range("F434").Value = 100112001
$ws.Range("G434").Value = "Berenjena"
$ws.Range("H434").Value = "Sin especificar"
$ws.Range("I434").Value = "Primera"
$ws.Range("J434").Value = 60
$ws.Range("K434").Value = 12000
$ws.Range("L434").Value = 12000
$ws.Range("M434").Value = 12000
$ws.Range("N434").Value = "$/caja 40 unidades"
$ws.Range("O434").Value = "Región de Arica y Parinacota"
$ws.Range("P434").Value = 300
$ws.Range("Q434").Value = 40
$ws.Range("R434").Value = "Hortaliza"
